$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit adds a new match (Welsh Premiership: The New Saints vs Bala Town)
# at row 10, shifting the previous row 10 (English Premier League) to row 11
# and the previous row 11 (Scottish Premiership) to row 12. It also refreshes
# a number of Back/Lay odds across the existing rows to their latest values.

# --- Insert new row at position 10, shifting rows 10-11 down to 11-12 ---
$ws.Rows.Item(10).Insert()

# --- New row 10: Welsh Premiership, The New Saints vs Bala Town ---
$ws.Range("A10").Value = "Welsh Premiership"
$ws.Range("B10:C10").NumberFormat = "@"
$ws.Range("B10").Value = "2026-01-06"
$ws.Range("C10").Value = "16:45:00"
$ws.Range("B10:C10").Style = "Normal"
$ws.Range("D10").Value = "The New Saints"
$ws.Range("E10").Value = "Bala Town"
$ws.Range("F10").Value = 1.09
$ws.Range("G10").Value = 1.23
$ws.Range("H10").Value = 15
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 8.6
$ws.Range("K10").Value = 950
$ws.Range("L10").Value = 1.01
$ws.Range("M10").Value = 1.01
$ws.Range("N10").Value = 1.78
$ws.Range("O10").Value = 1.1
$ws.Range("P10").Value = 1.78
$ws.Range("Q10").Value = 1.1
$ws.Range("R10").Value = 1.78
$ws.Range("S10").Value = 1.74
$ws.Range("T10").Value = 1.01
$ws.Range("U10").Value = 1.01
$ws.Range("V10").Value = 1.01
$ws.Range("W10").Value = 5.5
$ws.Range("X10").Value = 1000
$ws.Range("Y10").Value = 1000
$ws.Range("Z10").Value = 1000
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 1000
$ws.Range("AC10").Value = 1000
$ws.Range("AD10").Value = 1000
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 1000
$ws.Range("AG10").Value = 1000
$ws.Range("AH10").Value = 1000
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 1000
$ws.Range("AL10").Value = 1000
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 1000

# --- Update modified numeric cells in existing rows (2-9) ---
$ws.Range("J2").Value = 3.8
$ws.Range("M2").Value = 1.09
$ws.Range("N2").Value = 3.4
$ws.Range("P2").Value = 1.81
$ws.Range("T2").Value = 2.02
$ws.Range("AC2").Value = 8
$ws.Range("AJ2").Value = 19.5
$ws.Range("AK2").Value = 19.5

$ws.Range("F3").Value = 1.78
$ws.Range("G3").Value = 1.93
$ws.Range("H3").Value = 5.1
$ws.Range("I3").Value = 6.6
$ws.Range("J3").Value = 3.25
$ws.Range("N3").Value = 2.94
$ws.Range("O3").Value = 1.41
$ws.Range("Q3").Value = 2.24
$ws.Range("T3").Value = 2.04
$ws.Range("U3").Value = 1.79
$ws.Range("AC3").Value = 9.4

$ws.Range("G4").Value = 3.1
$ws.Range("H4").Value = 2.7
$ws.Range("T4").Value = 1.83
$ws.Range("U4").Value = 1.98
$ws.Range("W4").Value = 1.48

$ws.Range("F5").Value = 5.5
$ws.Range("G5").Value = 5.7
$ws.Range("H5").Value = 1.82
$ws.Range("I5").Value = 1.83
$ws.Range("O5").Value = 1.41
$ws.Range("P5").Value = 1.79
$ws.Range("Q5").Value = 2.22
$ws.Range("R5").Value = 1.3
$ws.Range("V5").Value = 2.2
$ws.Range("W5").Value = 1.21
$ws.Range("Y5").Value = 7.6
$ws.Range("AA5").Value = 18
$ws.Range("AF5").Value = 38
$ws.Range("AN5").Value = 110

$ws.Range("F6").Value = 2.5
$ws.Range("K6").Value = 4
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 2.16
$ws.Range("O6").Value = 1.04
$ws.Range("Q6").Value = 1.67
$ws.Range("T6").Value = 1.05

$ws.Range("F7").Value = 5.9
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 1.82
$ws.Range("I7").Value = 1.83
$ws.Range("J7").Value = 3.55
$ws.Range("K7").Value = 3.6
$ws.Range("N7").Value = 3
$ws.Range("Q7").Value = 2.44
$ws.Range("T7").Value = 2.22
$ws.Range("V7").Value = 2.2
$ws.Range("W7").Value = 1.2

$ws.Range("L8").Value = 1.43
$ws.Range("N8").Value = 3.7
$ws.Range("O8").Value = 1.34
$ws.Range("P8").Value = 1.94
$ws.Range("Q8").Value = 2.04
$ws.Range("R8").Value = 1.35
$ws.Range("T8").Value = 1.96
$ws.Range("U8").Value = 1.98
$ws.Range("Y8").Value = 8.199999999999999
$ws.Range("AB8").Value = 18
$ws.Range("AL8").Value = 85

$ws.Range("F9").Value = 2.92
$ws.Range("G9").Value = 2.96
$ws.Range("I9").Value = 2.64
$ws.Range("S9").Value = 3.85
$ws.Range("T9").Value = 1.84
$ws.Range("AB9").Value = 11.5
$ws.Range("AD9").Value = 12
$ws.Range("AF9").Value = 19
$ws.Range("AN9").Value = 36
$ws.Range("AO9").Value = 27

# --- Row 11 (shifted from old row 10): English Premier League, West Ham vs Nottm Forest ---
$ws.Range("F11").Value = 3.3
$ws.Range("G11").Value = 3.35
$ws.Range("L11").Value = 1.42
$ws.Range("P11").Value = 1.95
$ws.Range("T11").Value = 1.79
$ws.Range("Y11").Value = 10.5
$ws.Range("Z11").Value = 14
$ws.Range("AE11").Value = 25
$ws.Range("AG11").Value = 13.5
$ws.Range("AK11").Value = 38

# --- Row 12 (shifted from old row 11): Scottish Premiership, Rangers vs Aberdeen ---
$ws.Range("F12").Value = 1.56
$ws.Range("G12").Value = 1.58
$ws.Range("H12").Value = 7.2
$ws.Range("I12").Value = 7.6
$ws.Range("J12").Value = 4.4
$ws.Range("N12").Value = 5.1
$ws.Range("O12").Value = 1.22
$ws.Range("P12").Value = 2.28
$ws.Range("Q12").Value = 1.69
$ws.Range("R12").Value = 1.54
$ws.Range("S12").Value = 2.72
$ws.Range("T12").Value = 1.79
$ws.Range("U12").Value = 2.12
$ws.Range("V12").Value = 1.15
$ws.Range("W12").Value = 2.72
$ws.Range("X12").Value = 21
$ws.Range("Y12").Value = 28
$ws.Range("Z12").Value = 1000
$ws.Range("AC12").Value = 10
$ws.Range("AD12").Value = 28
$ws.Range("AE12").Value = 1000
$ws.Range("AF12").Value = 10
$ws.Range("AH12").Value = 21
$ws.Range("AI12").Value = 1000
$ws.Range("AJ12").Value = 14
$ws.Range("AL12").Value = 36
$ws.Range("AM12").Value = 120
$ws.Range("AN12").Value = 7
$ws.Range("AO12").Value = 980
